$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

# Fill in the Day 9 (row 54 header) totals that were previously left blank
$ws.Range("C55").Value = 718
$ws.Range("C56").Value = 903
$ws.Range("C57").Value = 636

# Update the active selection to reflect the last cell touched
$ws.Range("C57").Select()
